# Add 6000-element dataset results for Bubble Sort and Merge Sort,
# add chart/section title above the summary table, and update the
# active selection in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New title label spanning the merged cell C2:K2 (adds a new shared string).
$ws.Range("C2").Value = "Average Runtime per Number of Elements"

# --- Bubble Sort: 6000 elements (column L, rows 17-21; detail rows 18-21 were blank) ---
$ws.Range("L18").Value = 54.672600000000003
$ws.Range("L19").Value = 54.997700000000002
$ws.Range("L20").Value = 55.184199999999997
$ws.Range("L21").Value = 57.026800000000001
# L22 holds =AVERAGE(L17:L21) and recalculates automatically.

# --- Merge Sort: 6000 elements (column L, rows 24-28; all detail rows were blank) ---
$ws.Range("L24").Value = 5.4630000000000001
$ws.Range("L25").Value = 5.4660000000000002
$ws.Range("L26").Value = 5.6996000000000002
$ws.Range("L27").Value = 5.3846999999999996
$ws.Range("L28").Value = 5.4020999999999999
# L29 holds =AVERAGE(L24:L28) and recalculates automatically (no longer #DIV/0!).

# --- Propagate the newly completed 6000-element averages into the summary table ---
$ws.Range("E4").Value = 55.545540000000003
$ws.Range("E5").Value = 5.4830799999999993

# Update the sheet's active selection/view state.
[void]$ws.Range("I32").Select()
